$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant rows appended below the existing data (rows 3-4 -> now 3-8)
$data = @(
    @("Oribjonov Islombek Xusniddin o'g'li", "Mehnat muhofazasi va texnika xavfsizligi", "O'zbek tili", "Sirtqi", "AD5879006", "53006075130046", "Andijon viloyati", "Jalaquduq tuman", "998335567538"),
    @("Ismoilov Alisher Adxamovich", "Menejment", "O'zbek tili", "Sirtqi", "AD0716061", "51110056750049", "Toshkent shahri", "Yunusobod tumani", "998337034905"),
    @("Miryaxyoyev Miralisher Mirxaydar o'g'li", "Metrologiya va standartlashtirish", "O'zbek tili", "Sirtqi", "AD7947777", "30605986780015", "Toshkent viloyati", "Qibray tumani", "998937037330"),
    @("/start", "Menejment", "Rus tili", "Sirtqi", "AB6908896", "51503016520051", "Toshkent shahri", "Shayxontohur tumani", "998990677063")
)

$startRow = 5
# Columns holding digit-only strings (passport serial numbers, JSHIR,
# phone numbers) need to be pre-formatted as text; otherwise Excel's
# COM layer auto-converts a purely numeric string into a Number, which
# would drop meaning (and not match the source data's text formatting).
$textCols = @(6, 9)  # F = JSHIR, I = phone number

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    foreach ($col in $textCols) {
        $ws.Cells.Item($row, $col).NumberFormat = "@"
    }

    for ($j = 0; $j -lt $values.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
